$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Update the title (shared string) in A1: 13.07.2020 -> 14.07.2020
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Données COVID-19 Valais 14.07.2020"

# ---------------------------------------------------------------------------
# 2) Corrected historical values for rows 127-138 (dates 44013-44024)
# ---------------------------------------------------------------------------
$ws.Range("G127").Value = 7
$ws.Range("I127").Value = 1

$ws.Range("G128").Value = 9

$ws.Range("G129").Value = 9

$ws.Range("G130").Value = 7

$ws.Range("G131").Value = 7

$ws.Range("G132").Value = 9

$ws.Range("G133").Value = 7

$ws.Range("G134").Value = 7

$ws.Range("G135").Value = 8

$ws.Range("B136").Formula = "=B135+C136"
$ws.Range("C136").Value = 3
$ws.Range("E136").Value = 2
$ws.Range("G136").Value = 6
$ws.Range("I136").Value = 2

$ws.Range("C137").Value = 5
$ws.Range("E137").Value = 2
$ws.Range("G137").Value = 6

$ws.Range("E138").Value = 2
$ws.Range("G138").Value = 6

# ---------------------------------------------------------------------------
# 3) Row 140 did not exist before: first clone row139's current ("last row")
#    formatting down onto row140, before row139 itself is re-formatted as a
#    normal (non-last) row.
# ---------------------------------------------------------------------------
$ws.Range("A139:M139").Copy()
$ws.Range("A140:M140").PasteSpecial(-4122)

# Re-format row139 as a normal interior row (copy format from row138).
$ws.Range("A138:M138").Copy()
$ws.Range("A139:M139").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4) Fill in the data for 14.07.2020 (row139, date 44025) and the brand new
#    row for 15.07.2020 (row140, date 44026).
# ---------------------------------------------------------------------------
$ws.Range("A139").Value = 44025
$ws.Range("B139").Formula = "=B138+C139"
$ws.Range("C139").Value = 1
$ws.Range("D139").Value = 0
$ws.Range("E139").Value = 2
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 5
$ws.Range("H139").Formula = "=G139+E139"
$ws.Range("I139").Value = 1
$ws.Range("J139").Formula = "=J138+K139"
$ws.Range("K139").Formula = "=L139+M139"
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = 0

$ws.Range("A140").Value = 44026
$ws.Range("B140").Formula = "=B139+C140"
$ws.Range("C140").Value = 0
$ws.Range("D140").Value = 0
$ws.Range("E140").Value = 2
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 5
$ws.Range("H140").Formula = "=G140+E140"
$ws.Range("I140").Value = 0
$ws.Range("J140").Formula = "=J139+K140"
$ws.Range("K140").Formula = "=L140+M140"
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = 0

# ---------------------------------------------------------------------------
# 5) Update the view: select A1:M1 and scroll back to the top of the sheet.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A1:M1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
